$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.640.01"
$ws.Range("E2").Value = "'  +1.15%  "
$ws.Range("D3").Value = "'1.868.20"
$ws.Range("E3").Value = "'  +0.51%  "
$ws.Range("E4").Value = "'  +0.34%  "
$ws.Range("D5").Value = "'331.42"
$ws.Range("E5").Value = "'  +2.79%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "'  +0.28%  "
$ws.Range("D7").Value = "'0.4692"
$ws.Range("E7").Value = "'  +3.88%  "
$ws.Range("E8").Value = "'  +2.03%  "
$ws.Range("E9").Value = "'  -0.66%  "
$ws.Range("D10").Value = "'0.08057"
$ws.Range("E10").Value = "'  +2.14%  "
$ws.Range("E11").Value = "'  +0.08%  "
$ws.Range("D12").Value = "'21.79"
$ws.Range("E12").Value = "'  +1.93%  "
$ws.Range("D13").Value = "'1.847.51"
$ws.Range("E13").Value = "'  -0.74%  "
$ws.Range("D14").Value = "'5.940"
$ws.Range("E14").Value = "'  +1.12%  "
$ws.Range("D15").Value = "'7.132"
$ws.Range("E15").Value = "'  -0.37%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "'  +0.30%  "
$ws.Range("D17").Value = "'0.00001045"
$ws.Range("E17").Value = "'  +1.45%  "
$ws.Range("D18").Value = "'86.55"
$ws.Range("E18").Value = "'  +1.30%  "
$ws.Range("D19").Value = "'0.06612"
$ws.Range("E19").Value = "'  +1.35%  "
$ws.Range("D20").Value = "'17.18"
$ws.Range("E20").Value = "'  +0.79%  "
$ws.Range("E21").Value = "'  +0.35%  "
$ws.Range("D22").Value = "'27.657.16"
$ws.Range("E22").Value = "'  +1.20%  "
$ws.Range("D23").Value = "'5.490"
$ws.Range("E23").Value = "'  -0.41%  "
$ws.Range("D24").Value = "'10.97"
$ws.Range("E24").Value = "'  +1.90%  "
$ws.Range("D25").Value = "'2.311"
$ws.Range("E25").Value = "'  +1.85%  "
$ws.Range("D26").Value = "'2.087.16"
$ws.Range("E26").Value = "'  +0.28%  "
$ws.Range("D27").Value = "'158.91"
$ws.Range("E27").Value = "'  +4.67%  "
$ws.Range("D28").Value = "'20.18"
$ws.Range("E28").Value = "'  +2.34%  "
$ws.Range("D29").Value = "'2.089"
$ws.Range("E29").Value = "'  +1.34%  "
$ws.Range("D30").Value = "'5.545"
$ws.Range("E30").Value = "'  +0.72%  "
$ws.Range("D31").Value = "'122.05"
$ws.Range("E31").Value = "'  +1.38%  "
$ws.Range("D32").Value = "'0.9670"
$ws.Range("E32").Value = "'  +3.21%  "
$ws.Range("D33").Value = "'0.09481"
$ws.Range("E33").Value = "'  +1.96%  "
$ws.Range("D34").Value = "'1.441"
$ws.Range("E34").Value = "'  -2.34%  "
$ws.Range("E35").Value = "'  -0.05%  "
$ws.Range("D36").Value = "'5.314"
$ws.Range("E36").Value = "'  +0.73%  "
$ws.Range("D37").Value = "'0.02254"
$ws.Range("E37").Value = "'  +1.15%  "
$ws.Range("D38").Value = "'0.06080"
$ws.Range("E38").Value = "'  +1.55%  "
$ws.Range("D39").Value = "'1.224"
$ws.Range("E39").Value = "'  +1.01%  "
$ws.Range("D40").Value = "'8.122"
$ws.Range("E40").Value = "'  -1.81%  "
$ws.Range("E41").Value = "'  +0.29%  "
$ws.Range("E42").Value = "'  +1.17%  "
$ws.Range("D43").Value = "'0.1887"
$ws.Range("E43").Value = "'  -0.15%  "
$ws.Range("E44").Value = "'  +0.79%  "
$ws.Range("D45").Value = "'1.263"
$ws.Range("E45").Value = "'  +0.07%  "
$ws.Range("D46").Value = "'0.5696"
$ws.Range("E46").Value = "'  +1.20%  "
$ws.Range("D47").Value = "'12.25"
$ws.Range("E47").Value = "'  +2.12%  "
$ws.Range("D48").Value = "'3.388"
$ws.Range("E48").Value = "'  +0.97%  "
$ws.Range("D49").Value = "'1.931"
$ws.Range("E49").Value = "'  +0.35%  "
$ws.Range("D50").Value = "'0.06847"
$ws.Range("E50").Value = "'  +0.65%  "
$ws.Range("D51").Value = "'114.10"
$ws.Range("E51").Value = "'  +5.38%  "
